# Adding test cases in IAM module
# Fills in the "Results" column (D) for rows 38-45 on the "Test Cases"
# sheet with "SKIP", matching the pattern already used by the other rows
# in the same column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

for ($r = 38; $r -le 45; $r++) {
    $ws.Cells.Item($r, 4).Value = "SKIP"
}
